$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2331081081081081
$ws.Range("C2").Value = 0.4695945945945946
$ws.Range("J2").Value = 0.01689189189189189
$ws.Range("O2").Value = 0.003378378378378379
$ws.Range("P2").Value = 0.1554054054054054
$ws.Range("S2").Value = 0.1216216216216216
$ws.Range("B3").Value = 0.01428571428571429
$ws.Range("C3").Value = 0.03571428571428571
$ws.Range("J3").Value = 0.007142857142857143
$ws.Range("P3").Value = 0.6928571428571428
$ws.Range("S3").Value = 0.25
$ws.Range("J4").Value = 0.05263157894736842
$ws.Range("P4").Value = 0.6842105263157895
$ws.Range("S4").Value = 0.2631578947368421
$ws.Range("B6").Value = 0.05825242718446602
$ws.Range("D6").Value = 0.02427184466019417
$ws.Range("E6").Value = 0.004854368932038835
$ws.Range("F6").Value = 0.04854368932038835
$ws.Range("J6").Value = 0.2766990291262136
$ws.Range("O6").Value = 0.01941747572815534
$ws.Range("Q6").Value = 0.1262135922330097
$ws.Range("R6").Value = 0.07766990291262135
$ws.Range("S6").Value = 0.3640776699029126
$ws.Range("B7").Value = 0.075
$ws.Range("D7").Value = 0.0125
$ws.Range("F7").Value = 0.05833333333333333
$ws.Range("J7").Value = 0.125
$ws.Range("O7").Value = 0.025
$ws.Range("Q7").Value = 0.1833333333333333
$ws.Range("R7").Value = 0.08749999999999999
$ws.Range("S7").Value = 0.4333333333333333
$ws.Range("B8").Value = 0.0514018691588785
$ws.Range("D8").Value = 0.02102803738317757
$ws.Range("E8").Value = 0.002336448598130841
$ws.Range("F8").Value = 0.0630841121495327
$ws.Range("J8").Value = 0.08878504672897196
$ws.Range("O8").Value = 0.01869158878504673
$ws.Range("Q8").Value = 0.1985981308411215
$ws.Range("R8").Value = 0.09579439252336448
$ws.Range("S8").Value = 0.4602803738317757
$ws.Range("B9").Value = 0.1093117408906883
$ws.Range("D9").Value = 0.008097165991902834
$ws.Range("E9").Value = 0.004048582995951417
$ws.Range("F9").Value = 0.05263157894736842
$ws.Range("J9").Value = 0.1012145748987854
$ws.Range("O9").Value = 0.01214574898785425
$ws.Range("Q9").Value = 0.1821862348178138
$ws.Range("R9").Value = 0.06072874493927125
$ws.Range("S9").Value = 0.4696356275303644
$ws.Range("B10").Value = 0.1044558071585099
$ws.Range("D10").Value = 0.01460920379839299
$ws.Range("E10").Value = 0.0007304601899196494
$ws.Range("F10").Value = 0.0547845142439737
$ws.Range("J10").Value = 0.1183345507669832
$ws.Range("O10").Value = 0.01241782322863404
$ws.Range("Q10").Value = 0.2220598977355734
$ws.Range("R10").Value = 0.07669831994156319
$ws.Range("S10").Value = 0.39590942293645
$ws.Range("G11").Value = 0.1560846560846561
$ws.Range("J11").Value = 0.07936507936507936
$ws.Range("K11").Value = 0.208994708994709
$ws.Range("L11").Value = 0.5476190476190477
$ws.Range("S11").Value = 0.007936507936507936
$ws.Range("G12").Value = 0.7102803738317757
$ws.Range("J12").Value = 0.205607476635514
$ws.Range("K12").Value = 0.009345794392523364
$ws.Range("L12").Value = 0.02803738317757009
$ws.Range("S12").Value = 0.04672897196261682
$ws.Range("G13").Value = 0.7450980392156863
$ws.Range("J13").Value = 0.196078431372549
$ws.Range("S13").Value = 0.05882352941176471
$ws.Range("F15").Value = 0.015
$ws.Range("H15").Value = 0.115
$ws.Range("I15").Value = 0.08
$ws.Range("J15").Value = 0.35
$ws.Range("K15").Value = 0.1
$ws.Range("M15").Value = 0.005
$ws.Range("N15").Value = 0.005
$ws.Range("O15").Value = 0.015
$ws.Range("S15").Value = 0.315
$ws.Range("F16").Value = 0.01183431952662722
$ws.Range("H16").Value = 0.2071005917159763
$ws.Range("I16").Value = 0.136094674556213
$ws.Range("J16").Value = 0.3846153846153846
$ws.Range("K16").Value = 0.106508875739645
$ws.Range("M16").Value = 0.005917159763313609
$ws.Range("O16").Value = 0.01775147928994083
$ws.Range("S16").Value = 0.1301775147928994
$ws.Range("F17").Value = 0.02976190476190476
$ws.Range("H17").Value = 0.1527777777777778
$ws.Range("I17").Value = 0.07539682539682539
$ws.Range("J17").Value = 0.4384920634920635
$ws.Range("K17").Value = 0.1150793650793651
$ws.Range("M17").Value = 0.0119047619047619
$ws.Range("N17").Value = 0.003968253968253968
$ws.Range("O17").Value = 0.05158730158730158
$ws.Range("S17").Value = 0.121031746031746
$ws.Range("F18").Value = 0.01020408163265306
$ws.Range("H18").Value = 0.1836734693877551
$ws.Range("I18").Value = 0.1326530612244898
$ws.Range("J18").Value = 0.3571428571428572
$ws.Range("K18").Value = 0.09183673469387756
$ws.Range("M18").Value = 0.02551020408163265
$ws.Range("O18").Value = 0.03061224489795918
$ws.Range("S18").Value = 0.1683673469387755
$ws.Range("F19").Value = 0.01642710472279261
$ws.Range("H19").Value = 0.1765913757700205
$ws.Range("I19").Value = 0.09924709103353867
$ws.Range("J19").Value = 0.3805612594113621
$ws.Range("K19").Value = 0.1245722108145106
$ws.Range("M19").Value = 0.02600958247775496
$ws.Range("N19").Value = 0.001368925393566051
$ws.Range("O19").Value = 0.06297056810403832
$ws.Range("S19").Value = 0.1122518822724162
